$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column C: width (closest reachable value to the target 53.43 chars
# given this engine quantizes ColumnWidth to 1/6-character increments)
$ws.Columns.Item(3).ColumnWidth = 52.6666667

# Column C header + German translations.
# Cells are written in this specific order so the workbook's shared-string
# table is (re)built with the same ordering the source file has.
$ws.Range("C1").Value = "German"
$ws.Range("C2").Value = "Karma während der Krise"
$ws.Range("C3").Value = "Größerer Genuß beim nächsten Kaffee um Dein Café zu unterstützen!"
$ws.Range("C4").Value = "In Kurzform:"
$ws.Range("C5").Value = "Mit guten Gewissen den nächsten Kaffee in der Isolation trinken! Spende einfach (D)einen Beitrag für Deinen Coffeeshop nach Wahl! Wie? Einfach die Emailadresse angeben (findet sich auf der Webseite Deines Coffeeshops) und folge diese Paypalanweisungen! Genieße tolles Karma beim nächsten Kaffee!"
$ws.Range("C6").Value = "Gerade kleinere Unternehmen sind durch die Koronakrise betroffen und haben große Mühe diese zu überleben. In vielen Cafés und Restaurants ist der Betrieb in den nächsten Wochen sehr eingeschränkt, wenn nicht sogar ganz zum Stillstand gekommen. In den Niederlanden betrifft dies ca. 400.000 Angestellte.Mit dieser kleinen Geste können wir diese Unternehmen unterstützen und mit gutem Gewissen unseren Kaffee auch in der Isolation geniessen."
$ws.Range("C7").Value = "Untersuchungen haben ergeben, dass durch Unterstützung anderer Stress und Angstzustände gelindert und Glückshormone freigesetzt werden!"
$ws.Range("C14").Value = "Bitte teile diese Webseite indem Du sie weitergibst an alle Personen, die Du kennst um Dein Lieblingsort auf der ganzen Welt zu retten!"
$ws.Range("C15").Value = "Vielen Dank und bleib gesund!"
$ws.Range("C12").Value = "Dadurch genießt Du Deinen Kaffee zu Hause noch mehr und es macht Dich glücklich geholfen zu haben!"
$ws.Range("C10").Value = "Unter 'Wie sende ich eine Zahlung mit PayPal?' lernst Du wie es funktioniert!"
$ws.Range("C11").Value = "Und nächstes Mal mit Deinem Kaffee zu Hause anstelle in Deinem Coffeeshop spendest Du den Anteil, den Du sonst dort ausgibst, um Dein Café zu unterstützen!"
$ws.Range("C13").Value = "Übrigens, auch Deinen Lieblingsfriseur , Restaurant oder jedes andere Unternehmen kannst Du so in der Krise unterstützen!"
$ws.Range("C9").Value = "Wie funktioniert das? PayPal macht es einfach! Alles was benötigt wird ist die Emailadresse des Unternehmens, dass Du unterstützen möchtest. Sende einfach einen beliebigen Betrag an die Emailadresse, die Du auf deren Webseite findest."
$ws.Range("C8").Value = "Hast auch Du ein Lieblingsrestaurant, dass Du nach der Krise wieder besuchen willst? Vielleicht sind dies mittlerweile sogar Deine Freunde? Das nächste Mal, wenn Du Deinen Kaffee draussen genießen möchtest, denk dran und unterstütze Deinen Coffeeshop mit dem Betrag, den Du normalerweise dort ausgeben würdest!"

# Match the saved selection/view state from the edit.
[void]$ws.Range("C9").Select()
